# Commit: "Fruta / hortaliza, semanal"
# Weekly data refresh: a new observation is inserted at row 9 (pushing the
# previously-existing rows 9-51 down to 10-52), and its values are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9; everything below (old rows 9..51) shifts
# down to 10..52, carrying its existing values/styles with it.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with this week's observation.
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 44819
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 100112026
$ws.Range("G9").Value = "Haba"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 13000
$ws.Range("M9").Value = 12500
$ws.Range("N9").Value = "$/saco 25 kilos"
$ws.Range("O9").Value = "Provincia del Elquí"
$ws.Range("P9").Value = 500
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
